$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 08:47"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2727996
$ws.Range("C4").Value = 143
$ws.Range("D4").Value = 1143490
$ws.Range("E4").Value = 1454383
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 130123

# Ucrania (row 37)
$ws.Range("B37").Value = 44998
$ws.Range("C37").Value = 664
$ws.Range("D37").Value = 19548
$ws.Range("E37").Value = 24277
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 14
$ws.Range("H37").Value = 1173

# Afganistan now overtakes Suiza -> row 46 becomes Afganistan with fresh data,
# row 47 becomes Suiza keeping its previous snapshot values (rows swap order).
$ws.Range("A46").Value = "Afganistan"
$ws.Range("B46").Value = 31762
$ws.Range("C46").Value = 245
$ws.Range("D46").Value = 15661
$ws.Range("E46").Value = 15327
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 28
$ws.Range("H46").Value = 774

$ws.Range("A47").Value = "Suiza"
$ws.Range("B47").Value = 31714
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 29200
$ws.Range("E47").Value = 551
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 1963

# El Salvador (row 79)
$ws.Range("D79").Value = 3809
$ws.Range("E79").Value = 2447
$ws.Range("G79").Value = 8
$ws.Range("H79").Value = 182

# Georgia (row 138)
$ws.Range("B138").Value = 931
$ws.Range("C138").Value = 3
$ws.Range("E138").Value = 122
